# This script reorders a handful of data rows across the
# "CATEGORIA A - BASICOS", "CATEGORIA C - BAJO IMPACTO" and
# "CATEGORIA D - SIN VENTAS" sheets (per the upstream commit that fixed
# clasificacionabc.py's row ordering before the mailer ran).
#
# Strategy: for every affected block of rows, copy each source row's full
# contents (values + number formats + fill/interior formatting) into a
# scratch row far below the sheet's used range, then copy the scratch
# rows back into their new target positions, and finally clear the
# scratch rows so the sheet's used range/dimensions are unaffected.

$wb = $excel.ActiveWorkbook

function Permute-Rows {
    param(
        $ws,
        [int]$lastCol,
        [hashtable]$mapping,   # targetRow -> sourceRow (both refer to ORIGINAL/before content)
        [int]$scratchOffset = 1000
    )

    $lastColLetter = [string]([char](64 + $lastCol))

    # 1) Stash every distinct source row's current contents into a scratch row.
    $sources = $mapping.Values | Sort-Object -Unique
    foreach ($src in $sources) {
        $scratch = $src + $scratchOffset
        $srcRange = "A$src`:$lastColLetter$src"
        $scratchRange = "A$scratch`:$lastColLetter$scratch"
        $ws.Range($srcRange).Copy($ws.Range($scratchRange))
    }

    # 2) Write the scratch rows into their new target positions. Clear the
    #    target row first: Range.Copy only overwrites cells the source
    #    actually populated, so a target cell that used to hold a value but
    #    whose new source cell is blank would otherwise keep its stale value.
    foreach ($target in $mapping.Keys) {
        $source = $mapping[$target]
        $scratch = $source + $scratchOffset
        $scratchRange = "A$scratch`:$lastColLetter$scratch"
        $targetRange = "A$target`:$lastColLetter$target"
        $ws.Range($targetRange).Clear()
        $ws.Range($scratchRange).Copy($ws.Range($targetRange))
    }

    # 3) Clean up the scratch rows.
    foreach ($src in $sources) {
        $scratch = $src + $scratchOffset
        $scratchRange = "A$scratch`:$lastColLetter$scratch"
        $ws.Range($scratchRange).Clear()
    }
}

# ---------------------------------------------------------------------
# CATEGORIA A - BASICOS (sheet 1): rows 51 and 52 swap places.
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item(1)
Permute-Rows $wsA 23 @{ 51 = 52; 52 = 51 } 1000

# ---------------------------------------------------------------------
# CATEGORIA C - BAJO IMPACTO (sheet 3): two independent pair swaps.
# ---------------------------------------------------------------------
$wsC = $wb.Worksheets.Item(3)
Permute-Rows $wsC 23 @{ 16 = 17; 17 = 16 } 1000
Permute-Rows $wsC 23 @{ 25 = 26; 26 = 25 } 1000
Permute-Rows $wsC 23 @{ 46 = 47; 47 = 46 } 1000
Permute-Rows $wsC 23 @{ 63 = 64; 64 = 63 } 1000

# ---------------------------------------------------------------------
# CATEGORIA D - SIN VENTAS (sheet 4): rows 2-18 get reshuffled into a
# new order (full permutation, not just pairwise swaps).
# Mapping below reads as targetRow = sourceRow (content that ends up at
# targetRow originally lived at sourceRow).
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item(4)
$mappingD = @{
    2  = 18
    3  = 13
    4  = 11
    5  = 12
    6  = 5
    7  = 8
    8  = 16
    9  = 2
    10 = 7
    11 = 10
    12 = 17
    13 = 6
    14 = 4
    15 = 9
    16 = 3
    17 = 14
    18 = 15
}
Permute-Rows $wsD 23 $mappingD 1000
